$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 was stored as text "2" - convert it to a real number (2)
$ws.Range("B4").Value = 2

# Add new row 5 data (annotation row for Ruilin)
$ws.Range("A5").Value = "Ruilin"

# B5 keeps "3" as text (matches how the source data stores this column),
# so force text formatting before assigning to avoid Excel auto-converting
# the numeric-looking string into a number.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "3"

$ws.Range("C5").Value = "无"
$ws.Range("D5").Value = "QSN"
$ws.Range("E5").Value = "MET"
$ws.Range("F5").Value = "53dcf950-aee9-43ba-bb93-9e7c5cd5833d"
$ws.Range("G5").Value = "By5SY2gA-_annotated.xlsx"
$ws.Range("H5").Value = "For instance, what about averaging WordNet path-based distance metrics and distance in word embedding space (for word similarity), and other ways of applying the affect data to email tone prediction?"
